$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04917089449816103
$ws.Range("D2").Value = 0.0201798235502082
$ws.Range("E2").Value = 0.08907804039668576
$ws.Range("F2").Value = 0.973411283863939
$ws.Range("G2").Value = 0.8376653156757072
$ws.Range("H2").Value = 0.8417501268403811
$ws.Range("K2").Value = 1.18336423364957
$ws.Range("M2").Value = 0.3903609337213396
$ws.Range("N2").Value = 1.194751420916911

$ws.Range("C3").Value = 0.04365387833961165
$ws.Range("D3").Value = 0.02034135116610436
$ws.Range("E3").Value = 0.08409063378335446
$ws.Range("F3").Value = 0.9530165663767747
$ws.Range("G3").Value = 0.816356225585352
$ws.Range("H3").Value = 0.8378073356566631
$ws.Range("K3").Value = 1.046738008367527
$ws.Range("M3").Value = 0.3513415549217882
$ws.Range("N3").Value = 1.212704287276077

$ws.Range("C4").Value = 0.04028331649385564
$ws.Range("D4").Value = 0.02044642775869043
$ws.Range("E4").Value = 0.08110301514950891
$ws.Range("F4").Value = 0.9412631593784653
$ws.Range("G4").Value = 0.8040017444839123
$ws.Range("H4").Value = 0.835944170149304
$ws.Range("K4").Value = 0.9630244209873808
$ws.Range("M4").Value = 0.3275340707216046
$ws.Range("N4").Value = 1.224271010639573

$ws.Range("C5").Value = 0.03891395961441901
$ws.Range("D5").Value = 0.02049072960788223
$ws.Range("E5").Value = 0.07990406411087037
$ws.Range("F5").Value = 0.9366659079900899
$ws.Range("G5").Value = 0.7991493452811795
$ws.Range("H5").Value = 0.8353247065070519
$ws.Range("K5").Value = 0.9289542429186497
$ws.Range("M5").Value = 0.3178697135032209
$ws.Range("N5").Value = 1.229121011414338

$ws.Range("C6").Value = 0.03868682882986718
$ws.Range("D6").Value = 0.02049817540422971
$ws.Range("E6").Value = 0.07970609100909343
$ws.Range("F6").Value = 0.9359141260608794
$ws.Range("G6").Value = 0.7983545729736079
$ws.Range("H6").Value = 0.8352302749367198
$ws.Range("K6").Value = 0.9232995513605431
$ws.Range("M6").Value = 0.3162671986030077
$ws.Range("N6").Value = 1.229934586195914

$ws.Range("C7").Value = 0.04026483206604325
$ws.Range("D7").Value = 0.02044701922914527
$ws.Range("E7").Value = 0.08108677098459438
$ws.Range("F7").Value = 0.9412003817469525
$ws.Range("G7").Value = 0.8039355674704183
$ws.Range("H7").Value = 0.8359352504417075
$ws.Range("K7").Value = 0.9625647617746154
$ws.Range("M7").Value = 0.3274035830688788
$ws.Range("N7").Value = 1.224335867199005

$ws.Range("C8").Value = 0.04726506966014199
$ws.Range("D8").Value = 0.02023429357161177
$ws.Range("E8").Value = 0.08734273346382793
$ws.Range("F8").Value = 0.9662189469840996
$ws.Range("G8").Value = 0.8301657897390697
$ws.Range("H8").Value = 0.8402746194898185
$ws.Range("K8").Value = 1.136218843769029
$ws.Range("M8").Value = 0.3768754998494543
$ws.Range("N8").Value = 1.200828613440638

$ws.Range("C9").Value = 0.06113086700693771
$ws.Range("D9").Value = 0.01986398586745608
$ws.Range("E9").Value = 0.100215049990517
$ws.Range("F9").Value = 1.021432361485481
$ws.Range("G9").Value = 0.8874511918486689
$ws.Range("H9").Value = 0.8532316029470053
$ws.Range("K9").Value = 1.478179707266861
$ws.Range("M9").Value = 0.4751112752005042
$ws.Range("N9").Value = 1.159053477473734

$ws.Range("C10").Value = 0.07140923185518488
$ws.Range("D10").Value = 0.01962056527909439
$ws.Range("E10").Value = 0.1100583518448275
$ws.Range("F10").Value = 1.065823439068623
$ws.Range("G10").Value = 0.9331935180001096
$ws.Range("H10").Value = 0.8654957159989465
$ws.Range("K10").Value = 1.730360434750025
$ws.Range("M10").Value = 0.5480744317927417
$ws.Range("N10").Value = 1.131008367911803

$ws.Range("C11").Value = 0.0761064585615685
$ws.Range("D11").Value = 0.01951607268849109
$ws.Range("E11").Value = 0.1146240191649355
$ws.Range("F11").Value = 1.086865519702471
$ws.Range("G11").Value = 0.9548162504853792
$ws.Range("H11").Value = 0.8716783643879751
$ws.Range("K11").Value = 1.845306156872823
$ws.Range("M11").Value = 0.5814488436438268
$ws.Range("N11").Value = 1.11882736591819

$ws.Range("C12").Value = 0.07788837749667721
$ws.Range("D12").Value = 0.019477404414463
$ws.Range("E12").Value = 0.1163658552739335
$ws.Range("F12").Value = 1.094956851724902
$ws.Range("G12").Value = 0.9631228188649175
$ws.Range("H12").Value = 0.8741069408500266
$ws.Range("K12").Value = 1.888866704121085
$ws.Range("M12").Value = 0.5941138697596955
$ws.Range("N12").Value = 1.114297959110008

$ws.Range("C13").Value = 0.07750446656524446
$ws.Range("D13").Value = 0.01948569219773866
$ws.Range("E13").Value = 0.1159901407745849
$ws.Range("F13").Value = 1.093208742766421
$ws.Range("G13").Value = 0.9613285574438635
$ws.Range("H13").Value = 0.8735800102639359
$ws.Range("K13").Value = 1.879483682400007
$ws.Range("M13").Value = 0.5913850266049394
$ws.Range("N13").Value = 1.11526973903792

$ws.Range("C14").Value = 0.07625299387385098
$ws.Range("D14").Value = 0.01951287335256069
$ws.Range("E14").Value = 0.114767060597174
$ws.Range("F14").Value = 1.087528723612294
$ws.Range("G14").Value = 0.9554972536512594
$ws.Range("H14").Value = 0.8718764112699944
$ws.Range("K14").Value = 1.848889243675615
$ws.Range("M14").Value = 0.5824902615672158
$ws.Range("N14").Value = 1.118453058078227

$ws.Range("C15").Value = 0.07548684745893297
$ws.Range("D15").Value = 0.01952964001360868
$ws.Range("E15").Value = 0.1140195791166505
$ws.Range("F15").Value = 1.084065620599546
$ws.Range("G15").Value = 0.9519408871122721
$ws.Range("H15").Value = 0.8708442982556051
$ws.Range("K15").Value = 1.830153588766166
$ws.Range("M15").Value = 0.5770454750944936
$ws.Range("N15").Value = 1.120413787384134

$ws.Range("C16").Value = 0.07110269767534305
$ws.Range("D16").Value = 0.0196275199746303
$ws.Range("E16").Value = 0.1097617671419329
$ws.Range("F16").Value = 1.064465459658464
$ws.Range("G16").Value = 0.9317969281087528
$ws.Range("H16").Value = 0.8651038581166688
$ws.Range("K16").Value = 1.722853071095983
$ws.Range("M16").Value = 0.54589706517244
$ws.Range("N16").Value = 1.131816060321314

$ws.Range("C17").Value = 0.06841874107838919
$ws.Range("D17").Value = 0.01968916684191413
$ws.Range("E17").Value = 0.1071724469623163
$ws.Range("F17").Value = 1.052659502541673
$ws.Range("G17").Value = 0.9196488286962961
$ws.Range("H17").Value = 0.8617372725637154
$ws.Range("K17").Value = 1.657086203366532
$ws.Range("M17").Value = 0.5268357968760569
$ws.Range("N17").Value = 1.13895895998154

$ws.Range("C18").Value = 0.06687702334342305
$ws.Range("D18").Value = 0.01972521166237406
$ws.Range("E18").Value = 0.105691395304838
$ws.Range("F18").Value = 1.045948809672666
$ws.Range("G18").Value = 0.9127381257797822
$ws.Range("H18").Value = 0.8598576866583869
$ws.Range("K18").Value = 1.619280223794874
$ws.Range("M18").Value = 0.5158894620495573
$ws.Range("N18").Value = 1.143121643884079

$ws.Range("C19").Value = 0.06635536917357854
$ws.Range("D19").Value = 0.01973751656921507
$ws.Range("E19").Value = 0.105191346382064
$ws.Range("F19").Value = 1.043690351714261
$ws.Range("G19").Value = 0.9104113873913775
$ws.Range("H19").Value = 0.859231027731937
$ws.Range("K19").Value = 1.606483429800619
$ws.Range("M19").Value = 0.5121861597066442
$ws.Range("N19").Value = 1.144540369531726

$ws.Range("C20").Value = 0.06870424285764898
$ws.Range("D20").Value = 0.01968254363476163
$ws.Range("E20").Value = 0.1074472279655296
$ws.Range("F20").Value = 1.053908001405986
$ws.Range("G20").Value = 0.9209340805126089
$ws.Range("H20").Value = 0.8620897704430206
$ws.Range("K20").Value = 1.664084982848294
$ws.Range("M20").Value = 0.5288631172073224
$ws.Range("N20").Value = 1.138192966374558

$ws.Range("C21").Value = 0.07662049474195953
$ws.Range("D21").Value = 0.01950486511177729
$ws.Range("E21").Value = 0.1151259561081019
$ws.Range("F21").Value = 1.089193731296987
$ws.Range("G21").Value = 0.9572068213332727
$ws.Range("H21").Value = 0.8723744248369485
$ws.Range("K21").Value = 1.857874670474757
$ws.Range("M21").Value = 0.5851021364098585
$ws.Range("N21").Value = 1.117515777460693

$ws.Range("C22").Value = 0.08181280394936152
$ws.Range("D22").Value = 0.01939399433784317
$ws.Range("E22").Value = 0.120219869496168
$ws.Range("F22").Value = 1.112973369746982
$ws.Range("G22").Value = 0.9816045772585653
$ws.Range("H22").Value = 0.8796053896993783
$ws.Range("K22").Value = 1.984720546708729
$ws.Range("M22").Value = 0.6220144035981292
$ws.Range("N22").Value = 1.104487580705449

$ws.Range("C23").Value = 0.07903984311577972
$ws.Range("D23").Value = 0.01945268633611441
$ws.Range("E23").Value = 0.1174941587776246
$ws.Range("F23").Value = 1.100215611926899
$ws.Range("G23").Value = 0.9685193014881861
$ws.Range("H23").Value = 0.8756993026086661
$ws.Range("K23").Value = 1.917002710130532
$ws.Range("M23").Value = 0.6022991055833842
$ws.Range("N23").Value = 1.111396438786818

$ws.Range("C24").Value = 0.06857516338433811
$ws.Range("D24").Value = 0.01968553610788426
$ws.Range("E24").Value = 0.1073229758766701
$ws.Range("F24").Value = 1.053343316215731
$ws.Range("G24").Value = 0.9203527896376897
$ws.Range("H24").Value = 0.8619302320382189
$ws.Range("K24").Value = 1.660920821386185
$ws.Range("M24").Value = 0.527946527597976
$ws.Range("N24").Value = 1.13853909724548

$ws.Range("C25").Value = 0.0573642875761351
$ws.Range("D25").Value = 0.01995914122488607
$ws.Range("E25").Value = 0.09666623743375169
$ws.Range("F25").Value = 1.005829297462299
$ws.Range("G25").Value = 0.8713186745983421
$ws.Range("H25").Value = 0.8492468393872343
$ws.Range("K25").Value = 1.3855101466973
$ws.Range("M25").Value = 0.3903609337213396
$ws.Range("N25").Value = 1.169890796227889
